$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Gnas"
$ws.Range("C2").Value = "Gcgr"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 95.63567833333333
$ws.Range("H2").Value = 286.907035
$ws.Range("I2").Value = 0.2808828217467972
$ws.Range("J2").Value = 0.2808828217467972
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.07089166666666667
$ws.Range("N2").Value = 0.212675
$ws.Range("O2").Value = 0.07641114965052251
$ws.Range("P2").Value = 0.07641114965052251
$ws.Range("Q2").Value = 6.779772629847223
$ws.Range("R2").Value = 61.017953668625
$ws.Range("S2").Value = 0.02146257932675556
$ws.Range("T2").Value = 0.02146257932675556

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Gnas"
$ws.Range("C3").Value = "Gcgr"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 95.63567833333333
$ws.Range("H3").Value = 286.907035
$ws.Range("I3").Value = 0.2808828217467972
$ws.Range("J3").Value = 0.2808828217467972
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.8254096666666667
$ws.Range("N3").Value = 2.476229
$ws.Range("O3").Value = 0.8896744078427822
$ws.Range("P3").Value = 0.8896744078427823
$ws.Range("Q3").Value = 78.93861337455722
$ws.Range("R3").Value = 710.447520371015
$ws.Range("S3").Value = 0.2498942581107915
$ws.Range("T3").Value = 0.2498942581107916

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Gnas"
$ws.Range("C4").Value = "Gcgr"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 95.63567833333333
$ws.Range("H4").Value = 286.907035
$ws.Range("I4").Value = 0.2808828217467972
$ws.Range("J4").Value = 0.2808828217467972
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.03146466666666667
$ws.Range("N4").Value = 0.094394
$ws.Range("O4").Value = 0.03391444250669529
$ws.Range("P4").Value = 0.0339144425066953
$ws.Range("Q4").Value = 3.009144740198889
$ws.Range("R4").Value = 27.08230266179
$ws.Range("S4").Value = 0.009525984309250096
$ws.Range("T4").Value = 0.009525984309250099

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Gnas"
$ws.Range("C5").Value = "Gcgr"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 133.6830266666667
$ws.Range("H5").Value = 401.04908
$ws.Range("I5").Value = 0.392628215788982
$ws.Range("J5").Value = 0.392628215788982
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.07089166666666667
$ws.Range("N5").Value = 0.212675
$ws.Range("O5").Value = 0.07641114965052251
$ws.Range("P5").Value = 0.07641114965052251
$ws.Range("Q5").Value = 9.477012565444445
$ws.Range("R5").Value = 85.293113089
$ws.Range("S5").Value = 0.03000117335366955
$ws.Range("T5").Value = 0.03000117335366955

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Gnas"
$ws.Range("C6").Value = "Gcgr"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 133.6830266666667
$ws.Range("H6").Value = 401.04908
$ws.Range("I6").Value = 0.392628215788982
$ws.Range("J6").Value = 0.392628215788982
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.8254096666666667
$ws.Range("N6").Value = 2.476229
$ws.Range("O6").Value = 0.8896744078427822
$ws.Range("P6").Value = 0.8896744078427823
$ws.Range("Q6").Value = 110.3432624799245
$ws.Range("R6").Value = 993.08936231932
$ws.Range("S6").Value = 0.3493112753844306
$ws.Range("T6").Value = 0.3493112753844307

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Gnas"
$ws.Range("C7").Value = "Gcgr"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 133.6830266666667
$ws.Range("H7").Value = 401.04908
$ws.Range("I7").Value = 0.392628215788982
$ws.Range("J7").Value = 0.392628215788982
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.03146466666666667
$ws.Range("N7").Value = 0.094394
$ws.Range("O7").Value = 0.03391444250669529
$ws.Range("P7").Value = 0.0339144425066953
$ws.Range("Q7").Value = 4.206291873057778
$ws.Range("R7").Value = 37.85662685752001
$ws.Range("S7").Value = 0.01331576705088178
$ws.Range("T7").Value = 0.01331576705088179

# Row 8
$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Gnas"
$ws.Range("C8").Value = "Gcgr"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 111.16377
$ws.Range("H8").Value = 333.49131
$ws.Range("I8").Value = 0.3264889624642208
$ws.Range("J8").Value = 0.3264889624642208
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.07089166666666667
$ws.Range("N8").Value = 0.212675
$ws.Range("O8").Value = 0.07641114965052251
$ws.Range("P8").Value = 0.07641114965052251
$ws.Range("Q8").Value = 7.88058492825
$ws.Range("R8").Value = 70.92526435425
$ws.Range("S8").Value = 0.0249473969700974
$ws.Range("T8").Value = 0.0249473969700974

# Row 9
$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Gnas"
$ws.Range("C9").Value = "Gcgr"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 111.16377
$ws.Range("H9").Value = 333.49131
$ws.Range("I9").Value = 0.3264889624642208
$ws.Range("J9").Value = 0.3264889624642208
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.8254096666666667
$ws.Range("N9").Value = 2.476229
$ws.Range("O9").Value = 0.8896744078427822
$ws.Range("P9").Value = 0.8896744078427823
$ws.Range("Q9").Value = 91.75565034111
$ws.Range("R9").Value = 825.80085306999
$ws.Range("S9").Value = 0.29046887434756
$ws.Range("T9").Value = 0.29046887434756

# Row 10
$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Gnas"
$ws.Range("C10").Value = "Gcgr"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 111.16377
$ws.Range("H10").Value = 333.49131
$ws.Range("I10").Value = 0.3264889624642208
$ws.Range("J10").Value = 0.3264889624642208
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.03146466666666667
$ws.Range("N10").Value = 0.094394
$ws.Range("O10").Value = 0.03391444250669529
$ws.Range("P10").Value = 0.0339144425066953
$ws.Range("Q10").Value = 3.49773096846
$ws.Range("R10").Value = 31.47957871614
$ws.Range("S10").Value = 0.01107269114656341
$ws.Range("T10").Value = 0.01107269114656342
